$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Column G (status) width change: file width 17 -> 8 (ColumnWidth property has +0.83 offset)
$ws1.Columns.Item(7).ColumnWidth = 7.17

# Sheet1 (Players) changes
$ws1.Range("G3").Value = "Final"
$ws1.Range("G5").Value = "Final"
$ws1.Range("H5").Value = 12
$ws1.Range("I5").Value = 11
$ws1.Range("O6").Value = 13
$ws1.Range("G8").Value = "Final"
$ws1.Range("G12").Value = "Final"
$ws1.Range("G14").Value = "Final"
$ws1.Range("O14").Value = 23
$ws1.Range("D17").Value = "Keyshawn Hall"
$ws1.Range("E17").Value = "AUB"
$ws1.Range("F17").Value = "AUB@MISS"
$ws1.Range("G17").Value = "Final"
$ws1.Range("H17").Value = 20
$ws1.Range("I17").Value = 19
$ws1.Range("J17").Value = 6
$ws1.Range("K17").Value = 2
$ws1.Range("L17").Value = 1
$ws1.Range("M17").Value = 1
$ws1.Range("D18").Value = "Mark Mitchell"
$ws1.Range("E18").Value = "MIZ"
$ws1.Range("F18").Value = "UGA@MIZ"
$ws1.Range("H18").Value = 18
$ws1.Range("I18").Value = 18
$ws1.Range("J18").Value = 5
$ws1.Range("K18").Value = 3
$ws1.Range("L18").Value = 2
$ws1.Range("N18").Value = 2
$ws1.Range("O18").Value = 36
$ws1.Range("D19").Value = "Tyler Nickel"
$ws1.Range("E19").Value = "VAN"
$ws1.Range("F19").Value = "VAN@ARK"
$ws1.Range("G19").Value = "Final"
$ws1.Range("H19").Value = 14
$ws1.Range("I19").Value = 17
$ws1.Range("K19").Value = 0
$ws1.Range("M19").Value = 0
$ws1.Range("N19").Value = 0
$ws1.Range("O19").Value = 32
$ws1.Range("G20").Value = "Final"
$ws1.Range("O20").Value = 35
$ws1.Range("G23").Value = "Final"
$ws1.Range("O23").Value = 36
$ws1.Range("D27").Value = "Myles Stute"
$ws1.Range("E27").Value = "SC"
$ws1.Range("F27").Value = "OU@SC"
$ws1.Range("G27").Value = "Final"
$ws1.Range("H27").Value = 13
$ws1.Range("I27").Value = 11
$ws1.Range("J27").Value = 6
$ws1.Range("K27").Value = 1
$ws1.Range("L27").Value = 2
$ws1.Range("N27").Value = 2
$ws1.Range("O27").Value = 30
$ws1.Range("D28").Value = "Nick Pringle"
$ws1.Range("E28").Value = "ARK"
$ws1.Range("F28").Value = "VAN@ARK"
$ws1.Range("I28").Value = 8
$ws1.Range("J28").Value = 4
$ws1.Range("L28").Value = 0
$ws1.Range("N28").Value = 0
$ws1.Range("O28").Value = 16
$ws1.Range("D29").Value = "AJ Storr"
$ws1.Range("E29").Value = "MISS"
$ws1.Range("F29").Value = "AUB@MISS"
$ws1.Range("H29").Value = 12
$ws1.Range("I29").Value = 18
$ws1.Range("J29").Value = 5
$ws1.Range("K29").Value = 4
$ws1.Range("L29").Value = 1
$ws1.Range("N29").Value = 1
$ws1.Range("O29").Value = 33
$ws1.Range("G31").Value = "Final"
$ws1.Range("G32").Value = "Final"
$ws1.Range("O32").Value = 28
$ws1.Range("G38").Value = "Final"
$ws1.Range("O38").Value = 24
$ws1.Range("G41").Value = "Final"
$ws1.Range("G44").Value = "Final"
$ws1.Range("H44").Value = 19
$ws1.Range("I44").Value = 18
$ws1.Range("G52").Value = "Final"
$ws1.Range("H52").Value = 13
$ws1.Range("I52").Value = 17
$ws1.Range("O52").Value = 34
$ws1.Range("G53").Value = "Final"
$ws1.Range("I53").Value = 9
$ws1.Range("O53").Value = 28
$ws1.Range("O54").Value = 23
$ws1.Range("D55").Value = "Patton Pinkins"
$ws1.Range("E55").Value = "MISS"
$ws1.Range("F55").Value = "AUB@MISS"
$ws1.Range("H55").Value = 18
$ws1.Range("I55").Value = 15
$ws1.Range("K55").Value = 1
$ws1.Range("N55").Value = 0
$ws1.Range("O55").Value = 30
$ws1.Range("D56").Value = "Elijah Strong"
$ws1.Range("E56").Value = "SC"
$ws1.Range("F56").Value = "OU@SC"
$ws1.Range("I56").Value = 17
$ws1.Range("J56").Value = 5
$ws1.Range("K56").Value = 3
$ws1.Range("N56").Value = 1
$ws1.Range("O56").Value = 28
$ws1.Range("D57").Value = "Urban Klavzar"
$ws1.Range("E57").Value = "FLA"
$ws1.Range("F57").Value = "LSU@FLA"
$ws1.Range("G57").Value = "Final"
$ws1.Range("H57").Value = 16
$ws1.Range("I57").Value = 18
$ws1.Range("J57").Value = 3
$ws1.Range("K57").Value = 2
$ws1.Range("O57").Value = 24
$ws1.Range("G58").Value = "Final"
$ws1.Range("H58").Value = 15
$ws1.Range("K58").Value = 3
$ws1.Range("G59").Value = "Final"
$ws1.Range("G61").Value = "Final"
$ws1.Range("O61").Value = 36
$ws1.Range("G62").Value = "Final"
$ws1.Range("O62").Value = 23
$ws1.Range("G65").Value = "Final"
$ws1.Range("G66").Value = "Final"
$ws1.Range("G67").Value = "Final"
$ws1.Range("O67").Value = 10
$ws1.Range("G71").Value = "Final"
$ws1.Range("O71").Value = 14
$ws1.Range("G73").Value = "Final"
$ws1.Range("G74").Value = "Final"
$ws1.Range("G76").Value = "Final"
$ws1.Range("G79").Value = "Final"
$ws1.Range("D81").Value = "Blake Muschalek"
$ws1.Range("E81").Value = "AUB"
$ws1.Range("G81").Value = "Final"
$ws1.Range("H81").Value = 1
$ws1.Range("I81").Value = 0
$ws1.Range("J81").Value = 0
$ws1.Range("K81").Value = 1
$ws1.Range("L81").Value = 0
$ws1.Range("O81").Value = 6
$ws1.Range("D82").Value = "Hayden Assemian"
$ws1.Range("E82").Value = "SC"
$ws1.Range("F82").Value = "OU@SC"
$ws1.Range("G82").Value = "Final"
$ws1.Range("K82").Value = 0
$ws1.Range("M82").Value = 1
$ws1.Range("O82").Value = 3
$ws1.Range("D83").Value = "Kezza Giffa"
$ws1.Range("E83").Value = "MISS"
$ws1.Range("F83").Value = "AUB@MISS"
$ws1.Range("J83").Value = 1
$ws1.Range("K83").Value = 3
$ws1.Range("M83").Value = 0
$ws1.Range("O83").Value = 8
$ws1.Range("D84").Value = "Miles Keeffe"
$ws1.Range("E84").Value = "VAN"
$ws1.Range("F84").Value = "VAN@ARK"
$ws1.Range("G84").Value = "Final"
$ws1.Range("K84").Value = 0
$ws1.Range("L84").Value = 1
$ws1.Range("O84").Value = 1
$ws1.Range("D85").Value = "Zach Day"
$ws1.Range("E85").Value = "MISS"
$ws1.Range("F85").Value = "AUB@MISS"
$ws1.Range("I85").Value = 1
$ws1.Range("O85").Value = 6
$ws1.Range("G89").Value = "Final"
$ws1.Range("O89").Value = 8
$ws1.Range("G92").Value = "Final"
$ws1.Range("O92").Value = 9
$ws1.Range("G93").Value = "Final"
$ws1.Range("G99").Value = "Final"
$ws1.Range("G100").Value = "Final"
$ws1.Range("O100").Value = 6
$ws1.Range("G101").Value = "Final"

# Sheet2 (OwnerTotals) changes
$ws2.Range("A2").Value = "Clay"
$ws2.Range("B2").Value = 53
$ws2.Range("A3").Value = "Tar"
$ws2.Range("B3").Value = 52
